# See #284 - can now display game screen on Firefox and Opera.
# The on-board animated assets were superseded by new sprite-strip assets
# for the beach / cave / forest mini-games, and the "talking heads" assets
# were renamed to the new heads-*-strip.png naming scheme. A new row of
# tile offsets (row 122) was also added at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the on_board_*_animated_*.png rows (60-65) to the new
#     game-*-jump/rollover-strips.png sprite sheets ---
$ws.Range("A60").Value = "game-beach-jump-strips.png"
$ws.Range("A61").Value = "game-beach-rollover-strips.png"
$ws.Range("A62").Value = "game-cave-jump-strips.png"
$ws.Range("A63").Value = "game-cave-rollover-strips.png"
$ws.Range("A64").Value = "game-forest-jump-strips.png"
$ws.Range("A65").Value = "game-forest-rollover-strips.png"

# --- Rename the "talking heads" assets to the new heads-*-strip.png
#     naming scheme ---
$ws.Range("A117").Value = "collage/heads-eyes-strip.png"
$ws.Range("A115").Value = "heads-base.png"
$ws.Range("A116").Value = "heads-mouth-strip.png"
$ws.Range("D119").Value = "heads-right-eyes-strip.png"
$ws.Range("D120").Value = "heads-left-eyes-strip.png"

# --- Add the new row of tile-column offsets below the existing data ---
$ws.Range("A122").Value = 0
$ws.Range("B122").Value = 46
$ws.Range("C122").Value = 92
$ws.Range("D122").Value = 138
$ws.Range("E122").Value = 184
$ws.Range("F122").Value = 230
$ws.Range("G122").Value = 276
$ws.Range("H122").Value = 322
$ws.Range("I122").Value = 368
$ws.Range("J122").Value = 414
$ws.Range("K122").Value = 460
$ws.Range("L122").Value = 506

# --- Scroll the frozen pane down near the bottom of the sheet and leave
#     the selection on the last edited cell, matching where the author
#     left off editing ---
$excel.ActiveWindow.ScrollRow = 104
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D120").Select()
